# Refresh the crypto price/volume snapshot (cols D = Price, E = Volume(1h))
# with the latest scraped values. Cells that look like plain numbers
# (e.g. "0.998", "1.00") are forced to Text format first so Excel keeps
# them as literal strings instead of coercing them into floating point
# numbers (which would silently drop meaningful trailing zeros).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.097.46'
$ws.Range('D3').Value = '1.640.27'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.25%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.32'
$ws.Range('E5').Value = '  +0.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.503'
$ws.Range('E6').Value = '  +0.65%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.258'
$ws.Range('E8').Value = '  +0.52%  '
$ws.Range('E9').Value = '  +0.32%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.71'
$ws.Range('E10').Value = '  +0.65%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0790'
$ws.Range('E11').Value = '  -0.28%  '
$ws.Range('D12').Value = '1.868.44'
$ws.Range('E12').Value = '  +0.45%  '
$ws.Range('E13').Value = '  +0.43%  '
$ws.Range('D14').Value = '1.669.37'
$ws.Range('E14').Value = '  +2.09%  '
$ws.Range('E15').Value = '  -2.75%  '
$ws.Range('E16').Value = '  -0.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.25'
$ws.Range('E17').Value = '  +0.08%  '
$ws.Range('D18').Value = '26.106.84'
$ws.Range('E18').Value = '  +1.25%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '194.99'
$ws.Range('E20').Value = '  +1.35%  '
$ws.Range('E21').Value = '  -0.84%  '
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.79'
$ws.Range('E24').Value = '  -1.55%  '
$ws.Range('E25').Value = '  -0.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '142.60'
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  +0.86%  '
$ws.Range('E28').Value = '  +0.34%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.58'
$ws.Range('E29').Value = '  +0.58%  '
$ws.Range('E30').Value = '  +0.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0503'
$ws.Range('E31').Value = '  +2.22%  '
$ws.Range('E32').Value = '  +0.41%  '
$ws.Range('E33').Value = '  +0.23%  '
$ws.Range('E34').Value = '  +1.13%  '
$ws.Range('E35').Value = '  +1.46%  '
$ws.Range('E36').Value = '  +0.70%  '
$ws.Range('D37').Value = '1.134.00'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.552'
$ws.Range('E38').Value = '  +1.43%  '
$ws.Range('E39').Value = '  -1.11%  '
$ws.Range('E40').Value = '  +1.33%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.95'
$ws.Range('E42').Value = '  -0.81%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.48'
$ws.Range('E43').Value = '  -1.34%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.799'
$ws.Range('E44').Value = '  +0.13%  '
$ws.Range('D45').Value = '1.777.59'
$ws.Range('E45').Value = '  +0.48%  '
$ws.Range('D46').Value = '0.0₆0112'
$ws.Range('E46').Value = '  -0.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '56.69'
$ws.Range('E48').Value = '  +2.20%  '
$ws.Range('E49').Value = '  +3.71%  '
$ws.Range('E50').Value = '  +3.04%  '
$ws.Range('E51').Value = '  -0.09%  '
